# q2: systematic factor analysis (bivariate)
# Adds three new variables (env_perc, disp_inc, death_rate) in columns G:I,
# fills in previously-blank data points in columns B, D, G, H, I, removes
# the now-unused O and R placeholder columns, widens the M:N helper columns,
# and moves the active selection to N3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header labels (columns G, H, I) -----------------------------------
$ws.Range("G1").Value = "env_perc"
$ws.Range("H1").Value = "disp_inc"
$ws.Range("I1").Value = "death_rate"

# --- New data values for the three new columns ------------------------------
# Rows 2-4 and 5-7 had no entries at all yet in G/H/I; most of these land with
# the workbook's default (unstyled) format, except the few noted below that
# pick up the bold "stat" style already used by the rest of the table.

$ws.Range("H2").Value = 23112
$ws.Range("I2").Value = 495

$ws.Range("H3").Value = 23716
$ws.Range("I3").Value = 503

$ws.Range("H4").Value = 24296
$ws.Range("I4").Value = 494

$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 24999
$ws.Range("I5").Value = 480
$ws.Range("I5").Font.Bold = $true

$ws.Range("G6").Value = 2
$ws.Range("H6").Value = 25404
$ws.Range("I6").Value = 445
$ws.Range("I6").Font.Bold = $true

$ws.Range("G7").Value = 2
$ws.Range("G7").Font.Bold = $true
$ws.Range("H7").Value = 25978
$ws.Range("I7").Value = 422
$ws.Range("I7").Font.Bold = $true

# Rows 8-23 already carry the bold "stat" style on these cells (they were
# blank placeholders), so a plain value assignment keeps that formatting.

$ws.Range("D8").Value = 0
$ws.Range("G8").Value = 3
$ws.Range("H8").Value = 26088
$ws.Range("I8").Value = 395

$ws.Range("D9").Value = 2
$ws.Range("G9").Value = 3
$ws.Range("H9").Value = 26167
$ws.Range("I9").Value = 381

$ws.Range("D10").Value = 34
$ws.Range("G10").Value = 4
$ws.Range("H10").Value = 25460
$ws.Range("I10").Value = 372

$ws.Range("D11").Value = 100
$ws.Range("G11").Value = 5
$ws.Range("H11").Value = 24180
$ws.Range("I11").Value = 346

$ws.Range("G12").Value = 6
$ws.Range("H12").Value = 24758
$ws.Range("I12").Value = 336

$ws.Range("G13").Value = 7
$ws.Range("H13").Value = 24856
$ws.Range("I13").Value = 309

$ws.Range("G14").Value = 7
$ws.Range("H14").Value = 24711
$ws.Range("I14").Value = 305

$ws.Range("G15").Value = 7
$ws.Range("H15").Value = 24863
$ws.Range("I15").Value = 298

$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 25591
$ws.Range("I16").Value = 284

$ws.Range("G17").Value = 9
$ws.Range("H17").Value = 26098
$ws.Range("I17").Value = 290

$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 26392
$ws.Range("I18").Value = 276

$ws.Range("G19").Value = 9
$ws.Range("H19").Value = 27192
$ws.Range("I19").Value = 272

$ws.Range("G20").Value = 13
$ws.Range("H20").Value = 27334
$ws.Range("I20").Value = 266

$ws.Range("G21").Value = 25
$ws.Range("H21").Value = 28118
$ws.Range("I21").Value = 254

$ws.Range("B22").Value = 4184
$ws.Range("G22").Value = 25
$ws.Range("H22").Value = 23750
$ws.Range("I22").Value = 259

$ws.Range("B23").Value = 4849
$ws.Range("G23").Value = 30
$ws.Range("H23").Value = 26589
$ws.Range("I23").Value = 260

$ws.Range("G24").Value = 27

# --- Remove the now-unused helper columns O and R ---------------------------
$ws.Range("O8:O23").Clear()
$ws.Range("R8:R23").Clear()

# --- Widen the M:N helper columns slightly ----------------------------------
$ws.Range("M1:N1").ColumnWidth = 7.95

# --- Move the active selection ----------------------------------------------
$ws.Range("N3").Select()
